$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ "E" = 3; "G" = 5.451731666666667; "H" = 16.355195; "I" = 0.1922099906071488; "J" = 0.1922099906071488; "K" = 3; "M" = 8.131233999999999; "N" = 24.393702; "O" = 0.02090995573015822; "P" = 0.02090995573015823; "Q" = 44.32930588687667; "R" = 398.96375298189; "S" = 0.004019102394489609; "T" = 0.00401910239448961 }
    3 = @{ "E" = 3; "G" = 5.451731666666667; "H" = 16.355195; "I" = 0.1922099906071488; "J" = 0.1922099906071488; "K" = 3; "M" = 243.3763986666667; "N" = 730.1291960000001; "O" = 0.625857000534647; "P" = 0.6258570005346471; "Q" = 1326.822819530358; "R" = 11941.40537577322; "S" = 0.1202959681941828; "T" = 0.1202959681941829 }
    4 = @{ "E" = 3; "G" = 5.451731666666667; "H" = 16.355195; "I" = 0.1922099906071488; "J" = 0.1922099906071488; "K" = 3; "M" = 103.9426383333333; "N" = 311.827915; "O" = 0.2672947262403034; "P" = 0.2672947262403035; "Q" = 566.6673729187139; "R" = 5100.006356268426; "S" = 0.05137671681998914; "T" = 0.05137671681998915 }
    5 = @{ "E" = 3; "G" = 5.451731666666667; "H" = 16.355195; "I" = 0.1922099906071488; "J" = 0.1922099906071488; "K" = 3; "M" = 33.41874933333333; "N" = 100.256248; "O" = 0.08593831749489127; "P" = 0.08593831749489128; "Q" = 182.1900540009289; "R" = 1639.71048600836; "S" = 0.01651820319848723; "T" = 0.01651820319848723 }
    6 = @{ "E" = 3; "G" = 15.797976; "H" = 47.393928; "I" = 0.5569842765993244; "J" = 0.5569842765993244; "K" = 3; "M" = 8.131233999999999; "N" = 24.393702; "O" = 0.02090995573015822; "P" = 0.02090995573015823; "Q" = 128.457039582384; "R" = 1156.113356241456; "S" = 0.01164651656608607; "T" = 0.01164651656608608 }
    7 = @{ "E" = 3; "G" = 15.797976; "H" = 47.393928; "I" = 0.5569842765993244; "J" = 0.5569842765993244; "K" = 3; "M" = 243.3763986666667; "N" = 730.1291960000001; "O" = 0.625857000534647; "P" = 0.6258570005346471; "Q" = 3844.854505102433; "R" = 34603.69054592189; "S" = 0.3485925086974133; "T" = 0.3485925086974134 }
    8 = @{ "E" = 3; "G" = 15.797976; "H" = 47.393928; "I" = 0.5569842765993244; "J" = 0.5569842765993244; "K" = 3; "M" = 103.9426383333333; "N" = 311.827915; "O" = 0.2672947262403034; "P" = 0.2672947262403035; "Q" = 1642.08330576668; "R" = 14778.74975190012; "S" = 0.1488789597337699; "T" = 0.1488789597337699 }
    9 = @{ "E" = 3; "G" = 15.797976; "H" = 47.393928; "I" = 0.5569842765993244; "J" = 0.5569842765993244; "K" = 3; "M" = 33.41874933333333; "N" = 100.256248; "O" = 0.08593831749489127; "P" = 0.08593831749489128; "Q" = 527.9485999180159; "R" = 4751.537399262144; "S" = 0.04786629160205507; "T" = 0.04786629160205508 }
    10 = @{ "E" = 3; "G" = 3.587063; "H" = 10.761189; "I" = 0.1264679532473782; "J" = 0.1264679532473782; "K" = 3; "M" = 8.131233999999999; "N" = 24.393702; "O" = 0.02090995573015822; "P" = 0.02090995573015823; "Q" = 29.167248625742; "R" = 262.505237631678; "S" = 0.002644439303686397; "T" = 0.002644439303686398 }
    11 = @{ "E" = 3; "G" = 3.587063; "H" = 10.761189; "I" = 0.1264679532473782; "J" = 0.1264679532473782; "K" = 3; "M" = 243.3763986666667; "N" = 730.1291960000001; "O" = 0.625857000534647; "P" = 0.6258570005346471; "Q" = 873.0064747304494; "R" = 7857.058272574045; "S" = 0.07915085388316007; "T" = 0.07915085388316008 }
    12 = @{ "E" = 3; "G" = 3.587063; "H" = 10.761189; "I" = 0.1264679532473782; "J" = 0.1264679532473782; "K" = 3; "M" = 103.9426383333333; "N" = 311.827915; "O" = 0.2672947262403034; "P" = 0.2672947262403035; "Q" = 372.8487920878817; "R" = 3355.639128790935; "S" = 0.03380421694142944; "T" = 0.03380421694142945 }
    13 = @{ "E" = 3; "G" = 3.587063; "H" = 10.761189; "I" = 0.1264679532473782; "J" = 0.1264679532473782; "K" = 3; "M" = 33.41874933333333; "N" = 100.256248; "O" = 0.08593831749489127; "P" = 0.08593831749489128; "Q" = 119.8751592398747; "R" = 1078.876433158872; "S" = 0.01086844311910225; "T" = 0.01086844311910225 }
    14 = @{ "E" = 3; "G" = 3.526644000000001; "H" = 10.579932; "I" = 0.1243377795461487; "J" = 0.1243377795461487; "K" = 3; "M" = 8.131233999999999; "N" = 24.393702; "O" = 0.02090995573015822; "P" = 0.02090995573015823; "Q" = 28.675967598696; "R" = 258.083708388264; "S" = 0.002599897465896142; "T" = 0.002599897465896142 }
    15 = @{ "E" = 3; "G" = 3.526644000000001; "H" = 10.579932; "I" = 0.1243377795461487; "J" = 0.1243377795461487; "K" = 3; "M" = 243.3763986666667; "N" = 730.1291960000001; "O" = 0.625857000534647; "P" = 0.6258570005346471; "Q" = 858.3019160994082; "R" = 7724.717244894674; "S" = 0.0778176697598908; "T" = 0.07781766975989081 }
    16 = @{ "E" = 3; "G" = 3.526644000000001; "H" = 10.579932; "I" = 0.1243377795461487; "J" = 0.1243377795461487; "K" = 3; "M" = 103.9426383333333; "N" = 311.827915; "O" = 0.2672947262403034; "P" = 0.2672947262403035; "Q" = 366.56868182242; "R" = 3299.118136401781; "S" = 0.03323483274511501; "T" = 0.03323483274511502 }
    17 = @{ "E" = 3; "G" = 3.526644000000001; "H" = 10.579932; "I" = 0.1243377795461487; "J" = 0.1243377795461487; "K" = 3; "M" = 33.41874933333333; "N" = 100.256248; "O" = 0.08593831749489127; "P" = 0.08593831749489128; "Q" = 117.856031823904; "R" = 1060.704286415136; "S" = 0.01068537957524672; "T" = 0.01068537957524673 }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value = $data[$row][$col]
    }
}